# "did first event analysis"
# Adds a new "useful episodes" / comment column (N/O) of event-analysis
# notes to sheet 1 ("T1 in EDTA + CTZ") and sheet 2 ("T1 in zinc + CTZ"),
# extends sheet 2 with six new rows of "skipped episode ..." notes, and
# moves the active sheet / selections around to where the author left off.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # "T1 in EDTA + CTZ"
$ws2 = $wb.Worksheets.Item(2)   # "T1 in zinc + CTZ"

# ---------------------------------------------------------------------
# Sheet 1 ("T1 in EDTA + CTZ"): new column N notes next to M on rows 2-4
# ---------------------------------------------------------------------
$ws1.Range("N2").Value = "0-46"

$ws1.Range("N3").Value = "all?"
$ws1.Range("N3").Style = "Bad"

$ws1.Range("N4").Value = "344-466"

# ---------------------------------------------------------------------
# Sheet 2 ("T1 in zinc + CTZ"): new column N/O notes + new rows 13-19
# ---------------------------------------------------------------------

# Row 4: N4 picks up the grey "customFormat" look already used by the
# rest of that row (copy format from A4, which carries that style).
$ws2.Range("A4").Copy()
$ws2.Range("N4").PasteSpecial(-4122)   # xlPasteFormats
$ws2.Range("N4").Value = "all"
$excel.CutCopyMode = $false

# Row 9: N9 gets the grey format (copy from A9), O9 is plain text.
$ws2.Range("A9").Copy()
$ws2.Range("N9").PasteSpecial(-4122)   # xlPasteFormats
$ws2.Range("N9").Value = "1-243,![114,115,156,164-166,191,201,204,237]"
$excel.CutCopyMode = $false

$ws2.Range("O9").Value = "skipped episode 166 because"

# Row 10: just a new comment in O10
$ws2.Range("O10").Value = "skipped episode 329 because"

# Row 11: just a new comment in O11
$ws2.Range("O11").Value = "skipped episode 331 because"

# Row 12: O12 gets a new comment
$ws2.Range("O12").Value = "skipped episode 345 because"

# New rows 13-19: column O filled first (13 through 19), then column M
# (14 through 18) -- matches the order the notes were actually typed in.
$ws2.Range("O13").Value = "skipped episode 347 because"
$ws2.Range("O14").Value = "skipped episode 356 because"
$ws2.Range("O15").Value = "skipped episode 364 because"
$ws2.Range("O16").Value = "skipped episode 370 because"
$ws2.Range("O17").Value = "skipped episode 379 because"
$ws2.Range("O18").Value = "skipped episode 400 because"
$ws2.Range("O19").Value = "skipped episode 406 because"

$ws2.Range("M14").Value = "skipped episode 791 because no events detected"
$ws2.Range("M15").Value = "skipped episode 1224 because no events detected"
$ws2.Range("M16").Value = "skipped episode 1229 because no events detected"
$ws2.Range("M17").Value = "skipped episode 1231 because no events detected"
$ws2.Range("M18").Value = "skipped episode 1356 because no events detected"

# M9 and M12 are re-pointed at text already used elsewhere in the workbook.
$ws2.Range("M9").Value = "0.  , -0.57, -1.18, -1.75, -2.33"
$ws2.Range("M12").Value = "0. ,-0.59, -1.22, -1.82, -2.42"

# ---------------------------------------------------------------------
# View state: sheet 1 becomes the active tab, selections moved to where
# work left off on each sheet.
# ---------------------------------------------------------------------
$ws2.Range("N12").Select()
$ws1.Select()
$ws1.Range("D26").Select()
